$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 4999.5
$ws.Range("I10").Value = 4999
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 4999
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -4706
$ws.Range("N10").Value = -5586

$ws.Range("H13").Value = 18961.125
$ws.Range("J13").Value = 23464
$ws.Range("L13").Value = 23464
$ws.Range("N13").Value = -23802

$ws.Range("H51").Value = 3333.1667
$ws.Range("I51").Value = 3749.75
$ws.Range("K51").Value = 3749.75
$ws.Range("M51").Value = -3265.75

$ws.Range("H103").Value = 742.7143
$ws.Range("I103").Value = 975
$ws.Range("J103").Value = 649.8
$ws.Range("K103").Value = 2925
$ws.Range("L103").Value = 1949.4
$ws.Range("M103").Value = -2339
$ws.Range("N103").Value = -3121.4

$ws.Range("H137").Value = 11288.702
$ws.Range("I137").Value = 5708.727
$ws.Range("J137").Value = 18961.166
$ws.Range("K137").Value = 17126.181
$ws.Range("L137").Value = 56883.49800000001
$ws.Range("M137").Value = -14576.181
$ws.Range("N137").Value = -61983.49800000001

$ws.Range("H138").Value = 1383.6111
$ws.Range("I138").Value = 1254.25
$ws.Range("K138").Value = 3762.75
$ws.Range("M138").Value = 1377.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3199.1936
$ws.Range("J2").Value = 4279.5835
$ws.Range("L2").Value = 4279.5835
$ws.Range("N2").Value = -4505.5835

$ws.Range("H10").Value = 4622.75
$ws.Range("J10").Value = 1833.1666
$ws.Range("L10").Value = 1833.1666
$ws.Range("N10").Value = -2173.1666

$ws.Range("H11").Value = 11668.4
$ws.Range("J11").Value = 3050
$ws.Range("L11").Value = 3050
$ws.Range("N11").Value = -3338

$ws.Range("H12").Value = 5138.5
$ws.Range("J12").Value = 4375
$ws.Range("L12").Value = 4375
$ws.Range("N12").Value = -4721

$ws.Range("H29").Value = 16046.5
$ws.Range("J29").Value = 19992
$ws.Range("L29").Value = 19992
$ws.Range("N29").Value = -20608

$ws.Range("H35").Value = 5718.636
$ws.Range("I35").Value = 2672.8572
$ws.Range("J35").Value = 11048.75
$ws.Range("K35").Value = 2672.8572
$ws.Range("L35").Value = 11048.75
$ws.Range("M35").Value = -2266.8572
$ws.Range("N35").Value = -11860.75

$ws.Range("H74").Value = 9095.865
$ws.Range("I74").Value = 9779.656000000001
$ws.Range("K74").Value = 9779.656000000001
$ws.Range("M74").Value = -8905.656000000001

$ws.Range("H77").Value = 9095.865
$ws.Range("I77").Value = 9779.656000000001
$ws.Range("K77").Value = 48898.28000000001
$ws.Range("M77").Value = -44530.28000000001

$ws.Range("H116").Value = 3199.1936
$ws.Range("J116").Value = 4279.5835
$ws.Range("L116").Value = 4279.5835
$ws.Range("N116").Value = -8867.583500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3199.1936
$ws.Range("J3").Value = 4279.5835
$ws.Range("L3").Value = 4279.5835
$ws.Range("N3").Value = -4507.5835

$ws.Range("H31").Value = 3499.5
$ws.Range("I31").Value = 3999
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 3999
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -3747
$ws.Range("N31").Value = -3504

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 955
$ws.Range("I5").Value = 456.42856
$ws.Range("K5").Value = 456.42856
$ws.Range("M5").Value = -344.42856

$ws.Range("H10").Value = 2178.3333
$ws.Range("I10").Value = 581.2
$ws.Range("J10").Value = 4174.75
$ws.Range("K10").Value = 581.2
$ws.Range("L10").Value = 4174.75
$ws.Range("M10").Value = -442.2
$ws.Range("N10").Value = -4452.75

$ws.Range("H12").Value = 572.375
$ws.Range("J12").Value = 681.61536
$ws.Range("L12").Value = 681.61536
$ws.Range("N12").Value = -1021.61536

$ws.Range("H31").Value = 347190.62
$ws.Range("I31").Value = 61394.35
$ws.Range("J31").Value = 671093.0600000001
$ws.Range("K31").Value = 61394.35
$ws.Range("L31").Value = 671093.0600000001
$ws.Range("M31").Value = -61099.35
$ws.Range("N31").Value = -671683.0600000001

$ws.Range("H34").Value = 347190.62
$ws.Range("I34").Value = 61394.35
$ws.Range("J34").Value = 671093.0600000001
$ws.Range("K34").Value = 61394.35
$ws.Range("L34").Value = 671093.0600000001
$ws.Range("M34").Value = -61192.35
$ws.Range("N34").Value = -671497.0600000001

$ws.Range("H58").Value = 3861.8928
$ws.Range("I58").Value = 1958.8182
$ws.Range("K58").Value = 1958.8182
$ws.Range("M58").Value = -1755.8182

$ws.Range("H132").Value = 20451.291
$ws.Range("I132").Value = 13019.741
$ws.Range("K132").Value = 39059.223
$ws.Range("M132").Value = -36529.223

$ws.Range("H136").Value = 3861.8928
$ws.Range("I136").Value = 1958.8182
$ws.Range("K136").Value = 5876.4546
$ws.Range("M136").Value = -3326.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1984.8823
$ws.Range("I5").Value = 724.5
$ws.Range("J5").Value = 2152.9333
$ws.Range("K5").Value = 2173.5
$ws.Range("L5").Value = 6458.7999
$ws.Range("M5").Value = -2061.5
$ws.Range("N5").Value = -6682.7999

$ws.Range("H7").Value = 111508.2
$ws.Range("I7").Value = 495.33334
$ws.Range("K7").Value = 1486.00002
$ws.Range("M7").Value = -1374.00002

$ws.Range("H12").Value = 3416
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3416
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 10248
$ws.Range("N12").Value = -10594
$ws.Range("M12").ClearContents()

$ws.Range("H36").Value = 29246
$ws.Range("I36").Value = 176
$ws.Range("J36").Value = 40874
$ws.Range("K36").Value = 528
$ws.Range("L36").Value = 122622
$ws.Range("M36").Value = -359
$ws.Range("N36").Value = -122960

$ws.Range("H92").Value = 285.66666
$ws.Range("I92").Value = 174
$ws.Range("J92").Value = 317.57144
$ws.Range("K92").Value = 522
$ws.Range("L92").Value = 952.71432
$ws.Range("M92").Value = 726
$ws.Range("N92").Value = -3448.71432

$ws.Range("H135").Value = 1984.8823
$ws.Range("I135").Value = 724.5
$ws.Range("J135").Value = 2152.9333
$ws.Range("K135").Value = 6520.5
$ws.Range("L135").Value = 19376.3997
$ws.Range("M135").Value = -3985.5
$ws.Range("N135").Value = -24446.3997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5660
$ws.Range("I29").Value = 4500
$ws.Range("K29").Value = 4500
$ws.Range("M29").Value = -4210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 9999
$ws.Range("I26").Value = 9999
$ws.Range("K26").Value = 9999
$ws.Range("M26").Value = -9704

$ws.Range("H34").Value = 1766.6666
$ws.Range("I34").Value = 2150
$ws.Range("K34").Value = 2150
$ws.Range("M34").Value = -1978

$ws.Range("H131").Value = 523324
$ws.Range("J131").Value = 523324
$ws.Range("L131").Value = 523324
$ws.Range("N131").Value = -533404

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 775
$ws.Range("I7").Value = 575
$ws.Range("J7").Value = 975
$ws.Range("K7").Value = 575
$ws.Range("L7").Value = 975
$ws.Range("M7").Value = -462
$ws.Range("N7").Value = -1201

$ws.Range("H9").Value = 9858.111000000001
$ws.Range("I9").Value = 6960.5713
$ws.Range("J9").Value = 19999.5
$ws.Range("K9").Value = 6960.5713
$ws.Range("L9").Value = 19999.5
$ws.Range("M9").Value = -6820.5713
$ws.Range("N9").Value = -20279.5

$ws.Range("H28").Value = 8904.25
$ws.Range("J28").Value = 8904.25
$ws.Range("L28").Value = 8904.25
$ws.Range("N28").Value = -9600.25

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws.Range("H34").Value = 20012.5
$ws.Range("I34").Value = 20012.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 20012.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -19809.5
$ws.Range("N34").ClearContents()

$ws.Range("H96").Value = 995
$ws.Range("J96").Value = 990
$ws.Range("L96").Value = 990
$ws.Range("N96").Value = -3736

$ws.Range("H132").Value = 24501.023
$ws.Range("I132").Value = 22883.938
$ws.Range("K132").Value = 68651.814
$ws.Range("M132").Value = -66121.814
